$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.05089966666666667
$ws.Range("H2").Value = 0.152699
$ws.Range("I2").Value = 0.02671091810242436
$ws.Range("J2").Value = 0.03728162213961778
$ws.Range("M2").Value = 121.928739
$ws.Range("N2").Value = 365.786217
$ws.Range("O2").Value = 0.2282232151508951
$ws.Range("P2").Value = 0.2419720431319445
$ws.Range("Q2").Value = 6.206132172187
$ws.Range("R2").Value = 55.855189549683
$ws.Range("S2").Value = 0.006096051608967534
$ws.Range("T2").Value = 0.009021110280396453
$ws.Range("G3").Value = 0.05089966666666667
$ws.Range("H3").Value = 0.152699
$ws.Range("I3").Value = 0.02671091810242436
$ws.Range("J3").Value = 0.03728162213961778
$ws.Range("O3").Value = 0.2768624053389947
$ws.Range("P3").Value = 0.2935413991166814
$ws.Range("Q3").Value = 7.528790092223333
$ws.Range("R3").Value = 67.75911083001
$ws.Range("S3").Value = 0.007395249034650105
$ws.Range("T3").Value = 0.01094369952420285
$ws.Range("G4").Value = 0.05089966666666667
$ws.Range("H4").Value = 0.152699
$ws.Range("I4").Value = 0.02671091810242436
$ws.Range("J4").Value = 0.03728162213961778
$ws.Range("M4").Value = 83.50496933333334
$ws.Range("N4").Value = 250.514908
$ws.Range("O4").Value = 0.1563025480180701
$ws.Range("P4").Value = 0.1657186665504434
$ws.Range("Q4").Value = 4.250375104076889
$ws.Range("R4").Value = 38.253375936692
$ws.Range("S4").Value = 0.004174984559310921
$ws.Range("T4").Value = 0.006178260707814947
$ws.Range("G5").Value = 0.05089966666666667
$ws.Range("H5").Value = 0.152699
$ws.Range("I5").Value = 0.02671091810242436
$ws.Range("J5").Value = 0.03728162213961778
$ws.Range("M5").Value = 91.06846250000001
$ws.Range("N5").Value = 182.136925
$ws.Range("O5").Value = 0.1704597085236707
$ws.Range("P5").Value = 0.1204857969594293
$ws.Range("Q5").Value = 4.635354385095834
$ws.Range("R5").Value = 27.812126310575
$ws.Range("S5").Value = 0.004553135314138897
$ws.Range("T5").Value = 0.004491905955432154
$ws.Range("G6").Value = 0.05089966666666667
$ws.Range("H6").Value = 0.152699
$ws.Range("I6").Value = 0.02671091810242436
$ws.Range("J6").Value = 0.03728162213961778
$ws.Range("M6").Value = 89.83562999999999
$ws.Range("N6").Value = 269.50689
$ws.Range("O6").Value = 0.1681521229683693
$ws.Range("P6").Value = 0.1782820942415013
$ws.Range("Q6").Value = 4.57260362179
$ws.Range("R6").Value = 41.15343259611
$ws.Range("S6").Value = 0.004491497585356903
$ws.Range("T6").Value = 0.006646645671771379
$ws.Range("I7").Value = 0.1226793098007496
$ws.Range("J7").Value = 0.1712289953794413
$ws.Range("M7").Value = 121.928739
$ws.Range("N7").Value = 365.786217
$ws.Range("O7").Value = 0.2282232151508951
$ws.Range("P7").Value = 0.2419720431319445
$ws.Range("Q7").Value = 28.503850316812
$ws.Range("R7").Value = 256.534652851308
$ws.Range("S7").Value = 0.02799826651521979
$ws.Range("T7").Value = 0.0414326298553937
$ws.Range("I8").Value = 0.1226793098007496
$ws.Range("J8").Value = 0.1712289953794413
$ws.Range("O8").Value = 0.2768624053389947
$ws.Range("P8").Value = 0.2935413991166814
$ws.Range("Q8").Value = 34.57862319097333
$ws.Range("S8").Value = 0.03396528879676324
$ws.Range("T8").Value = 0.05026279887302497
$ws.Range("I9").Value = 0.1226793098007496
$ws.Range("J9").Value = 0.1712289953794413
$ws.Range("M9").Value = 83.50496933333334
$ws.Range("N9").Value = 250.514908
$ws.Range("O9").Value = 0.1563025480180701
$ws.Range("P9").Value = 0.1657186665504434
$ws.Range("Q9").Value = 19.52134637091022
$ws.Range("R9").Value = 175.692117338192
$ws.Range("S9").Value = 0.01917508871095536
$ws.Range("T9").Value = 0.02837584078905304
$ws.Range("I10").Value = 0.1226793098007496
$ws.Range("J10").Value = 0.1712289953794413
$ws.Range("M10").Value = 91.06846250000001
$ws.Range("N10").Value = 182.136925
$ws.Range("O10").Value = 0.1704597085236707
$ws.Range("P10").Value = 0.1204857969594293
$ws.Range("Q10").Value = 21.28949946478334
$ws.Range("R10").Value = 127.7369967887
$ws.Range("S10").Value = 0.02091187939052088
$ws.Range("T10").Value = 0.02063066197085443
$ws.Range("I11").Value = 0.1226793098007496
$ws.Range("J11").Value = 0.1712289953794413
$ws.Range("M11").Value = 89.83562999999999
$ws.Range("N11").Value = 269.50689
$ws.Range("O11").Value = 0.1681521229683693
$ws.Range("P11").Value = 0.1782820942415013
$ws.Range("Q11").Value = 21.00129445804
$ws.Range("R11").Value = 189.01165012236
$ws.Range("S11").Value = 0.02062878638729032
$ws.Range("T11").Value = 0.03052706389111514
$ws.Range("G12").Value = 1.620901
$ws.Range("H12").Value = 3.241802
$ws.Range("I12").Value = 0.8506097720968261
$ws.Range("J12").Value = 0.791489382480941
$ws.Range("M12").Value = 121.928739
$ws.Range("N12").Value = 365.786217
$ws.Range("O12").Value = 0.2282232151508951
$ws.Range("P12").Value = 0.2419720431319445
$ws.Range("Q12").Value = 197.634414973839
$ws.Range("R12").Value = 1185.806489843034
$ws.Range("S12").Value = 0.1941288970267078
$ws.Range("T12").Value = 0.1915183029961544
$ws.Range("G13").Value = 1.620901
$ws.Range("H13").Value = 3.241802
$ws.Range("I13").Value = 0.8506097720968261
$ws.Range("J13").Value = 0.791489382480941
$ws.Range("O13").Value = 0.2768624053389947
$ws.Range("P13").Value = 0.2935413991166814
$ws.Range("Q13").Value = 239.75448541133
$ws.Range("R13").Value = 1438.52691246798
$ws.Range("S13").Value = 0.2355018675075814
$ws.Range("T13").Value = 0.2323349007194536
$ws.Range("G14").Value = 1.620901
$ws.Range("H14").Value = 3.241802
$ws.Range("I14").Value = 0.8506097720968261
$ws.Range("J14").Value = 0.791489382480941
$ws.Range("M14").Value = 83.50496933333334
$ws.Range("N14").Value = 250.514908
$ws.Range("O14").Value = 0.1563025480180701
$ws.Range("P14").Value = 0.1657186665504434
$ws.Range("Q14").Value = 135.3532882973693
$ws.Range("R14").Value = 812.1197297842159
$ws.Range("S14").Value = 0.1329524747478038
$ws.Range("T14").Value = 0.1311645650535754
$ws.Range("G15").Value = 1.620901
$ws.Range("H15").Value = 3.241802
$ws.Range("I15").Value = 0.8506097720968261
$ws.Range("J15").Value = 0.791489382480941
$ws.Range("M15").Value = 91.06846250000001
$ws.Range("N15").Value = 182.136925
$ws.Range("O15").Value = 0.1704597085236707
$ws.Range("P15").Value = 0.1204857969594293
$ws.Range("Q15").Value = 147.6129619347125
$ws.Range("R15").Value = 590.45184773885
$ws.Range("S15").Value = 0.144994693819011
$ws.Range("T15").Value = 0.09536322903314276
$ws.Range("G16").Value = 1.620901
$ws.Range("H16").Value = 3.241802
$ws.Range("I16").Value = 0.8506097720968261
$ws.Range("J16").Value = 0.791489382480941
$ws.Range("M16").Value = 89.83562999999999
$ws.Range("N16").Value = 269.50689
$ws.Range("O16").Value = 0.1681521229683693
$ws.Range("P16").Value = 0.1782820942415013
$ws.Range("Q16").Value = 145.61466250263
$ws.Range("R16").Value = 873.6879750157799
$ws.Range("S16").Value = 0.1430318389957221
$ws.Range("T16").Value = 0.1411083846786148
